# Commit: Sat, Jul 04, 2020  7:07:25 AM
#
# The table on slide 5 (the "B1 - TYPES OF FINANCIAL DOCUMENTS" slide) had
# its table style switched from the deck's custom "Table_0" style
# ({3CA4073A-4532-42D6-A11F-646A34C93FA0}) to the built-in PowerPoint table
# style {ECB27656-76DE-41F2-ACA5-255B3BF364EB} (one of the gallery's
# "Medium Style" table styles), which is what happens when someone picks a
# different style from the Table Styles gallery on the Table Design tab.

$p = $ppt.ActivePresentation

# Slide 5 -> the slide with the financial-documents table; shape 2 on that
# slide is the graphicFrame that holds the table.
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)

if ($sh.HasTable) {
    $tbl = $sh.Table
    $tbl.ApplyStyle("{ECB27656-76DE-41F2-ACA5-255B3BF364EB}")
}
